# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest counts pulled from the source site.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 16513
$ws1.Range("F5").Value  = 744
$ws1.Range("F6").Value  = 15634
$ws1.Range("F9").Value  = 27
$ws1.Range("F11").Value = 652
$ws1.Range("F15").Value = 1165
$ws1.Range("F19").Value = 558
$ws1.Range("F20").Value = 45
$ws1.Range("F21").Value = 47
$ws1.Range("F23").Value = 92
$ws1.Range("F25").Value = 74
$ws1.Range("F26").Value = 278
$ws1.Range("F27").Value = 382
$ws1.Range("F30").Value = 5807
$ws1.Range("F31").Value = 5265

# Sheet "全部类型" (Worksheets index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 16513
$ws4.Range("F5").Value  = 744
$ws4.Range("F6").Value  = 15634
$ws4.Range("F9").Value  = 27
$ws4.Range("F11").Value = 652
$ws4.Range("F15").Value = 1165
$ws4.Range("F19").Value = 558
$ws4.Range("F20").Value = 45
$ws4.Range("F21").Value = 47
$ws4.Range("F25").Value = 92
$ws4.Range("F27").Value = 74
$ws4.Range("F28").Value = 278
$ws4.Range("F29").Value = 382
$ws4.Range("F32").Value = 5807
$ws4.Range("F34").Value = 5265
